# Actualización automática 2025-06-02 13:21:56
# Adds a new "PRESUPUESTO" (budget) column G to the "VENTA MENSUAL" sheet,
# mirroring the formatting of the existing "junio" column F, and fills in
# the per-client budget figures plus the column total in row 22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Clone the formatting (header style, currency number format, totals style,
# borders, etc.) of column F into the new column G before writing values.
$ws.Range("F1:F22").Copy()
$ws.Range("G1:G22").PasteSpecial(-4122)

# Header
$ws.Range("G1").Value = "PRESUPUESTO"

# Per-row budget values (rows 2-21 = clients, row 22 = column total)
$valores = @(0, 1000, 0, 3000, 5000, 0, 6000, 0, 7000, 0, 6000, 6000, 1000, 400, 6500, 0, 4000, 0, 500, 4000, 50400)

for ($i = 0; $i -lt $valores.Length; $i++) {
    $fila = $i + 2
    $ws.Cells.Item($fila, 7).Value = $valores[$i]
}

# Match column G's width (17) to the rest of the sheet's custom widths.
$ws.Columns("G").ColumnWidth = 16.17
